$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "呼唤效果发动时，如果发动该效果的牌由玩家控制，则玩家可以选择将被呼唤的牌放置在其左侧或右侧。"
$ws.Range("A11").Value = "战旗"
$ws.Range("D11").Value = "Warbanner"
$ws.Range("B11").Value = 1

$ws.Range("B11").WrapText = $true
$ws.Range("C11").WrapText = $true
Write-Host "WrapText set"
Write-Host "Row11 height auto:" $ws.Rows.Item(11).RowHeight

Write-Host "Done"
